# Progress.xlsx update — append Run 2 generations 130-137 (program had not
# converged; continuing the run produced 8 more generations of data) and
# refresh the chart series ranges to include them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Run 2")

# ---------------------------------------------------------------------
# 1. Append the new data rows (A = BestIndividual_Training_Error stays
#    flat at the same value already used for rows 122-129; B = new
#    Average_Training_Error readings; C = %Change formula).
# ---------------------------------------------------------------------
$newB = @{
    130 = 30181015.068997201
    131 = 25263459.829275001
    132 = 25418286.691030201
    133 = 25564998.2441779
    134 = 27015850.369339801
    135 = 27399332.242027398
    136 = 32855212.3060726
    137 = 28522028.6611492
}

for ($r = 130; $r -le 137; $r++) {
    $ws.Cells.Item($r, 1).Value = 12.328828005938
    $ws.Cells.Item($r, 1).NumberFormat = "#,##0.00"

    $ws.Cells.Item($r, 2).Value = $newB[$r]
    $ws.Cells.Item($r, 2).NumberFormat = "#,##0.00"
}

# Percent-change column — enter as two fill operations (130:135 then
# 136:137) matching how the formula was originally extended, each one
# becomes its own shared-formula group.
$ws.Range("C130:C135").Formula = "=(B130-B129)/B129"
$ws.Range("C136:C137").Formula = "=(B136-B135)/B135"

# ---------------------------------------------------------------------
# 2. Update the sheet view: scrolled one column back and the current
#    selection moved to C25.
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("C25").Select()

# ---------------------------------------------------------------------
# 3. Extend both chart series (Average_Training_Error / column B, and
#    BestIndividual_Training_Error / column A) to cover the new rows.
# ---------------------------------------------------------------------
$chartObj = $ws.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection()

$sB = $series.Item(1)
$sB.Formula = "=SERIES('Run 2'!`$B`$1,,'Run 2'!`$B`$2:`$B`$137,1)"

$sA = $series.Item(2)
$sA.Formula = "=SERIES('Run 2'!`$A`$1,,'Run 2'!`$A`$2:`$A`$137,2)"

"Run 2 extended through row 137"
